# Actualización automática 2025-07-03 13:15:08
# Updates the "CUMPLIMIENTO MENSUAL" sheet (row 2 = OTROS, row 4 = TOTAL)
# with refreshed VENTA / POR CUMPLIR / CUMPLIMIENTO figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 - "OTROS" group
$ws.Range("D2").Value = -42.54
$ws.Range("E2").Value = 42.54

# Row 4 - "TOTAL" row
$ws.Range("D4").Value = 289.21
$ws.Range("E4").Value = 13434.13
$ws.Range("F4").Value = 0.02107431572780387
